$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2564746666666666
$ws.Range("H2").Value = 0.7694239999999999
$ws.Range("I2").Value = 0.1818007399394835
$ws.Range("J2").Value = 0.1818007399394835
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.04738633333333334
$ws.Range("N2").Value = 0.142159
$ws.Range("O2").Value = 0.05760194168856402
$ws.Range("P2").Value = 0.05760194168856402
$ws.Range("Q2").Value = 0.01215339404622222
$ws.Range("R2").Value = 0.109380546416
$ws.Range("S2").Value = 0.01047207562093192
$ws.Range("T2").Value = 0.01047207562093192
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2564746666666666
$ws.Range("H3").Value = 0.7694239999999999
$ws.Range("I3").Value = 0.1818007399394835
$ws.Range("J3").Value = 0.1818007399394835
$ws.Range("O3").Value = 0.7659981644722047
$ws.Range("P3").Value = 0.7659981644722047
$ws.Range("Q3").Value = 0.1616174257084444
$ws.Range("R3").Value = 1.454556831376
$ws.Range("S3").Value = 0.139259033093333
$ws.Range("T3").Value = 0.139259033093333
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2564746666666666
$ws.Range("H4").Value = 0.7694239999999999
$ws.Range("I4").Value = 0.1818007399394835
$ws.Range("J4").Value = 0.1818007399394835
$ws.Range("M4").Value = 0.1451156666666667
$ws.Range("N4").Value = 0.435347
$ws.Range("O4").Value = 0.1763998938392313
$ws.Range("P4").Value = 0.1763998938392313
$ws.Range("Q4").Value = 0.03721849223644444
$ws.Range("R4").Value = 0.334966430128
$ws.Range("S4").Value = 0.03206963122521859
$ws.Range("T4").Value = 0.03206963122521859
$ws.Range("I5").Value = 0.7694380609030022
$ws.Range("J5").Value = 0.7694380609030022
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.04738633333333334
$ws.Range("N5").Value = 0.142159
$ws.Range("O5").Value = 0.05760194168856402
$ws.Range("P5").Value = 0.05760194168856402
$ws.Range("Q5").Value = 0.05143699608388889
$ws.Range("R5").Value = 0.462932964755
$ws.Range("S5").Value = 0.0443211263170965
$ws.Range("T5").Value = 0.0443211263170965
$ws.Range("I6").Value = 0.7694380609030022
$ws.Range("J6").Value = 0.7694380609030022
$ws.Range("O6").Value = 0.7659981644722047
$ws.Range("P6").Value = 0.7659981644722047
$ws.Range("S6").Value = 0.5893881423267521
$ws.Range("T6").Value = 0.5893881423267521
$ws.Range("I7").Value = 0.7694380609030022
$ws.Range("J7").Value = 0.7694380609030022
$ws.Range("M7").Value = 0.1451156666666667
$ws.Range("N7").Value = 0.435347
$ws.Range("O7").Value = 0.1763998938392313
$ws.Range("P7").Value = 0.1763998938392313
$ws.Range("S7").Value = 0.1357287922591536
$ws.Range("T7").Value = 0.1357287922591536
$ws.Range("G8").Value = 0.06878966666666667
$ws.Range("I8").Value = 0.0487611991575143
$ws.Range("J8").Value = 0.0487611991575143
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.04738633333333334
$ws.Range("N8").Value = 0.142159
$ws.Range("O8").Value = 0.05760194168856402
$ws.Range("P8").Value = 0.05760194168856402
$ws.Range("Q8").Value = 0.003259690074555556
$ws.Range("R8").Value = 0.029337210671
$ws.Range("S8").Value = 0.002808739750535596
$ws.Range("T8").Value = 0.002808739750535596
$ws.Range("G9").Value = 0.06878966666666667
$ws.Range("I9").Value = 0.0487611991575143
$ws.Range("J9").Value = 0.0487611991575143
$ws.Range("O9").Value = 0.7659981644722047
$ws.Range("P9").Value = 0.7659981644722047
$ws.Range("Q9").Value = 0.04334778552011111
$ws.Range("S9").Value = 0.03735098905211957
$ws.Range("T9").Value = 0.03735098905211957
$ws.Range("G10").Value = 0.06878966666666667
$ws.Range("I10").Value = 0.0487611991575143
$ws.Range("J10").Value = 0.0487611991575143
$ws.Range("M10").Value = 0.1451156666666667
$ws.Range("N10").Value = 0.435347
$ws.Range("O10").Value = 0.1763998938392313
$ws.Range("P10").Value = 0.1763998938392313
$ws.Range("Q10").Value = 0.009982458338111112
$ws.Range("R10").Value = 0.089842125043
$ws.Range("S10").Value = 0.008601470354859137
$ws.Range("T10").Value = 0.008601470354859137
